$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) and Column E (Volume 1h) updates.
# D-column values are forced to text via a leading apostrophe (Excel's
# standard 'store as text' convention) since several of them are
# numeric-looking strings (e.g. '566.23') that Excel would otherwise
# auto-convert to a number; Style is reset to Normal afterwards so no
# stray quote-prefix formatting is left behind on the cell.
$ws.Range("D2").Value = "'63.157.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.70%  "
$ws.Range("D3").Value = "'2.589.79"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.36%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'566.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.92%  "
$ws.Range("D6").Value = "'152.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.73%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'0.613"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.43%  "
$ws.Range("D9").Value = "'2.590.39"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.17%  "
$ws.Range("E10").Value = "  -8.56%  "
$ws.Range("D11").Value = "'5.73"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("D13").Value = "'0.373"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.94%  "
$ws.Range("D14").Value = "'27.67"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.44%  "
$ws.Range("D15").Value = "'3.065.03"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.98%  "
$ws.Range("D16").Value = "'0.0000175"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -8.96%  "
$ws.Range("D17").Value = "'63.105.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.60%  "
$ws.Range("D18").Value = "'2.620.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.40%  "
$ws.Range("D19").Value = "'11.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.23%  "
$ws.Range("D20").Value = "'7.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.45%  "
$ws.Range("D21").Value = "'4.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.83%  "
$ws.Range("D22").Value = "'336.80"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.21%  "
$ws.Range("D24").Value = "'66.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.77%  "
$ws.Range("D25").Value = "'1.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("E26").Value = "  -7.39%  "
$ws.Range("D27").Value = "'8.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.59%  "
$ws.Range("D28").Value = "'571.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.60%  "
$ws.Range("D29").Value = "'1.51"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.78%  "
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("D31").Value = "'0.158"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.73%  "
$ws.Range("D32").Value = "'7.63"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.10%  "
$ws.Range("D33").Value = "'2.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.36%  "
$ws.Range("D34").Value = "'1.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.89%  "
$ws.Range("D35").Value = "'6.41"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.39%  "
$ws.Range("D36").Value = "'5.27"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.15%  "
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("D38").Value = "'0.395"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.00%  "
$ws.Range("D39").Value = "'19.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.83%  "
$ws.Range("E40").Value = "  +0.89%  "
$ws.Range("D41").Value = "'1.83"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.27%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").Value = "'41.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.64%  "
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").Value = "'155.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.95%  "
$ws.Range("D46").Value = "'22.94"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("D47").Value = "'3.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.64%  "
$ws.Range("D48").Value = "'0.0573"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.26%  "
$ws.Range("D49").Value = "'0.623"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.97%  "
$ws.Range("D50").Value = "'0.0988"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.58%  "
$ws.Range("D51").Value = "'0.0242"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.74%  "
